$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin name / link text cells (rows 17-24 shifted)
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

# Update numeric-looking text cells (Price / Volume) - force text storage
$deCells = @(
    @{Cell="D2"; Value="321.30"},
    @{Cell="E2"; Value="6.03%"},
    @{Cell="D3"; Value="49.16"},
    @{Cell="E3"; Value="11.53%"},
    @{Cell="D4"; Value="5.342"},
    @{Cell="E4"; Value="5.09%"},
    @{Cell="D5"; Value="0.08053"},
    @{Cell="E5"; Value="4.58%"},
    @{Cell="D6"; Value="4.607"},
    @{Cell="E6"; Value="4.20%"},
    @{Cell="D7"; Value="1.364"},
    @{Cell="E7"; Value="30.56%"},
    @{Cell="D8"; Value="1.643"},
    @{Cell="E8"; Value="1.46%"},
    @{Cell="D9"; Value="0.1278"},
    @{Cell="E9"; Value="-0.80%"},
    @{Cell="E10"; Value="5.28%"},
    @{Cell="D11"; Value="0.09659"},
    @{Cell="E11"; Value="4.14%"},
    @{Cell="D12"; Value="0.04719"},
    @{Cell="E12"; Value="12.55%"},
    @{Cell="D13"; Value="0.1047"},
    @{Cell="E13"; Value="-0.03%"},
    @{Cell="D14"; Value="0.001324"},
    @{Cell="E14"; Value="3.29%"},
    @{Cell="D15"; Value="0.04202"},
    @{Cell="E15"; Value="0.65%"},
    @{Cell="D16"; Value="0.005864"},
    @{Cell="E16"; Value="1.69%"},
    @{Cell="D17"; Value="3.347"},
    @{Cell="E17"; Value="0.08%"},
    @{Cell="D18"; Value="2.441"},
    @{Cell="E18"; Value="4.74%"},
    @{Cell="D19"; Value="0.3506"},
    @{Cell="E19"; Value="4.67%"},
    @{Cell="D20"; Value="8.015"},
    @{Cell="E20"; Value="-7.34%"},
    @{Cell="D21"; Value="0.1368"},
    @{Cell="E21"; Value="-2.28%"},
    @{Cell="D22"; Value="0.3093"},
    @{Cell="E22"; Value="-2.66%"},
    @{Cell="D23"; Value="0.001316"},
    @{Cell="E23"; Value="2.41%"},
    @{Cell="D24"; Value="0.004333"},
    @{Cell="E24"; Value="-1.89%"},
    @{Cell="E25"; Value="-0.03%"},
    @{Cell="D26"; Value="0.0003545"},
    @{Cell="E26"; Value="-95.27%"},
    @{Cell="D38"; Value="0.02731"},
    @{Cell="E38"; Value="9.80%"},
    @{Cell="D39"; Value="0.06062"},
    @{Cell="E39"; Value="14.34%"},
    @{Cell="D40"; Value="0.01085"},
    @{Cell="D41"; Value="0.008021"},
    @{Cell="E41"; Value="3.46%"},
    @{Cell="D42"; Value="0.1463"},
    @{Cell="E42"; Value="8.59%"},
    @{Cell="D43"; Value="0.007911"},
    @{Cell="E43"; Value="7.48%"},
    @{Cell="D44"; Value="0.008667"},
    @{Cell="E44"; Value="14.91%"},
    @{Cell="D45"; Value="0.3494"},
    @{Cell="E45"; Value="15.83%"},
    @{Cell="D46"; Value="0.00006860"},
    @{Cell="E46"; Value="2.92%"},
    @{Cell="D47"; Value="0.00000000751"},
    @{Cell="E47"; Value="0.24%"},
    @{Cell="D48"; Value="0.05876"},
    @{Cell="E48"; Value="37.90%"},
    @{Cell="D49"; Value="0.004007"},
    @{Cell="E49"; Value="-4.61%"},
    @{Cell="D50"; Value="0.00002104"},
    @{Cell="E50"; Value="0.24%"},
    @{Cell="D51"; Value="0.0002004"},
    @{Cell="E51"; Value="0.24%"}
)
foreach ($item in $deCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}
